$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.759.61"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "1.546.43"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'205.88"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.246"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'21.37"
$ws.Range("E9").Value = "  -4.06%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "'0.0854"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "1.769.38"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "1.554.28"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "'3.66"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "'0.510"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "26.762.21"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "'61.09"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "'213.21"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'8.96"
$ws.Range("E23").Value = "  -4.65%  "
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").Value = "'153.10"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "'14.89"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").Value = "'0.0460"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "1.345.73"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("D34").Value = "'2.90"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "'1.50"
$ws.Range("E35").Value = "  -3.91%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "'0.925"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").Value = "'0.518"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").Value = "'0.799"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.67"
$ws.Range("E42").Value = "  +5.04%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'0.992"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.19"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'1.75"
$ws.Range("E45").Value = "  -4.30%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'62.78"
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.682.46"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D48").Value = "'2.25"
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'85.64"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0512"
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0973"
$ws.Range("E51").Value = "  -1.91%  "
